$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 86.75
$ws.Range("I2").Value = 83.5
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 83.5
$ws.Range("L2").Value = 90
$ws.Range("M2").Value = 29.5
$ws.Range("N2").Value = -316
$ws.Range("H17").Value = 1456
$ws.Range("J17").Value = 1456
$ws.Range("L17").Value = 4368
$ws.Range("N17").Value = -4704
$ws.Range("H40").Value = 1052.6316
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1000
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -1350
$ws.Range("H43").Value = 1460.1923
$ws.Range("J43").Value = 1603.0476
$ws.Range("L43").Value = 1603.0476
$ws.Range("N43").Value = -1741.0476
$ws.Range("H100").Value = 1596
$ws.Range("I100").Value = 1216.3572
$ws.Range("J100").Value = 2659
$ws.Range("K100").Value = 1216.3572
$ws.Range("L100").Value = 2659
$ws.Range("M100").Value = -675.3571999999999
$ws.Range("N100").Value = -3741
$ws.Range("H137").Value = 3706112.8
$ws.Range("I137").Value = 9092542
$ws.Range("J137").Value = 2942.9375
$ws.Range("K137").Value = 27277626
$ws.Range("L137").Value = 8828.8125
$ws.Range("M137").Value = -27275076
$ws.Range("N137").Value = -13928.8125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1665.7715
$ws.Range("I2").Value = 1806.1428
$ws.Range("J2").Value = 1455.2142
$ws.Range("K2").Value = 1806.1428
$ws.Range("L2").Value = 1455.2142
$ws.Range("M2").Value = -1693.1428
$ws.Range("N2").Value = -1681.2142
$ws.Range("H61").Value = 30365478
$ws.Range("I61").Value = 35751450
$ws.Range("J61").Value = 204054.8
$ws.Range("K61").Value = 35751450
$ws.Range("L61").Value = 204054.8
$ws.Range("M61").Value = -35751238
$ws.Range("N61").Value = -204478.8
$ws.Range("H116").Value = 1665.7715
$ws.Range("I116").Value = 1806.1428
$ws.Range("J116").Value = 1455.2142
$ws.Range("K116").Value = 1806.1428
$ws.Range("L116").Value = 1455.2142
$ws.Range("M116").Value = 487.8571999999999
$ws.Range("N116").Value = -6043.2142
$ws.Range("H136").Value = 30365478
$ws.Range("I136").Value = 35751450
$ws.Range("J136").Value = 204054.8
$ws.Range("K136").Value = 107254350
$ws.Range("L136").Value = 612164.3999999999
$ws.Range("M136").Value = -107251800
$ws.Range("N136").Value = -617264.3999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1665.7715
$ws.Range("I3").Value = 1806.1428
$ws.Range("J3").Value = 1455.2142
$ws.Range("K3").Value = 1806.1428
$ws.Range("L3").Value = 1455.2142
$ws.Range("M3").Value = -1692.1428
$ws.Range("N3").Value = -1683.2142
$ws.Range("H134").Value = 3748.6924
$ws.Range("I134").Value = 3748.6924
$ws.Range("K134").Value = 11246.0772
$ws.Range("M134").Value = -8711.0772

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 136937.5
$ws.Range("I134").Value = 2500
$ws.Range("J134").Value = 181750
$ws.Range("K134").Value = 7500
$ws.Range("L134").Value = 545250
$ws.Range("M134").Value = -4965
$ws.Range("N134").Value = -550320

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 650.46344
$ws.Range("I107").Value = 606.48
$ws.Range("J107").Value = 719.1875
$ws.Range("K107").Value = 1819.44
$ws.Range("L107").Value = 2157.5625
$ws.Range("M107").Value = 100.5599999999999
$ws.Range("N107").Value = -5997.5625
$ws.Range("H113").Value = 635.0968
$ws.Range("I113").Value = 551.6875
$ws.Range("J113").Value = 724.06665
$ws.Range("K113").Value = 1655.0625
$ws.Range("L113").Value = 2172.19995
$ws.Range("M113").Value = 514.9375
$ws.Range("N113").Value = -6512.19995
$ws.Range("H131").Value = 934.5
$ws.Range("I131").Value = 396.84616
$ws.Range("J131").Value = 1302.3684
$ws.Range("K131").Value = 1190.53848
$ws.Range("L131").Value = 3907.1052
$ws.Range("M131").Value = 3849.46152
$ws.Range("N131").Value = -13987.1052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 54029.973
$ws.Range("I132").Value = 41182.88
$ws.Range("J132").Value = 78735.92
$ws.Range("K132").Value = 123548.64
$ws.Range("L132").Value = 236207.76
$ws.Range("M132").Value = -121018.64
$ws.Range("N132").Value = -241267.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 948.38464
$ws.Range("I16").Value = 939.0909
$ws.Range("K16").Value = 939.0909
$ws.Range("M16").Value = -769.0909
$ws.Range("H35").Value = 2337
$ws.Range("I35").Value = 2337
$ws.Range("K35").Value = 2337
$ws.Range("M35").Value = -2001
$ws.Range("H40").Value = 2881.077
$ws.Range("I40").Value = 2595.4
$ws.Range("J40").Value = 3833.3333
$ws.Range("K40").Value = 2595.4
$ws.Range("L40").Value = 3833.3333
$ws.Range("M40").Value = -2459.4
$ws.Range("N40").Value = -4105.3333
$ws.Range("H82").Value = 2622.8572
$ws.Range("I82").Value = 5000
$ws.Range("J82").Value = 2226.6667
$ws.Range("K82").Value = 5000
$ws.Range("L82").Value = 2226.6667
$ws.Range("M82").Value = -4639
$ws.Range("N82").Value = -2948.6667
$ws.Range("H85").Value = 2622.8572
$ws.Range("I85").Value = 5000
$ws.Range("J85").Value = 2226.6667
$ws.Range("K85").Value = 5000
$ws.Range("L85").Value = 2226.6667
$ws.Range("M85").Value = -3752
$ws.Range("N85").Value = -4722.6667
$ws.Range("H100").Value = 1532.9565
$ws.Range("I100").Value = 1311.2858
$ws.Range("J100").Value = 1877.7778
$ws.Range("K100").Value = 1311.2858
$ws.Range("L100").Value = 1877.7778
$ws.Range("M100").Value = -770.2858000000001
$ws.Range("N100").Value = -2959.7778
$ws.Range("H122").Value = 3356.1702
$ws.Range("I122").Value = 2941.4285
$ws.Range("J122").Value = 3532.121
$ws.Range("K122").Value = 8824.2855
$ws.Range("L122").Value = 10596.363
$ws.Range("M122").Value = -6374.2855
$ws.Range("N122").Value = -15496.363

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2006.091
$ws.Range("I81").Value = 1302.8182
$ws.Range("J81").Value = 2709.3635
$ws.Range("K81").Value = 2605.6364
$ws.Range("L81").Value = 5418.727
$ws.Range("M81").Value = -1544.6364
$ws.Range("N81").Value = -7540.727
$ws.Range("H84").Value = 2006.091
$ws.Range("I84").Value = 1302.8182
$ws.Range("J84").Value = 2709.3635
$ws.Range("K84").Value = 13028.182
$ws.Range("L84").Value = 27093.635
$ws.Range("M84").Value = -7724.181999999999
$ws.Range("N84").Value = -37701.63499999999
$ws.Range("H96").Value = 1474.2307
$ws.Range("I96").Value = 1423.4445
$ws.Range("J96").Value = 1588.5
$ws.Range("K96").Value = 1423.4445
$ws.Range("L96").Value = 1588.5
$ws.Range("M96").Value = -50.44450000000006
$ws.Range("N96").Value = -4334.5
